$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$studentIds = @{
    2  = "16-0062"
    3  = "17-0108"
    4  = "16-0087"
    5  = "16-0044"
    6  = "16-0166"
    7  = "18-0098"
    8  = "16-0101"
    9  = "16-0075"
    10 = "16-0028"
    11 = "18-0120"
    12 = "18-0232"
    13 = "16-0130"
    14 = "18-0175"
    15 = "16-0073"
}

foreach ($row in 2..15) {
    $src = $ws.Range("G$row")
    $dst = $ws.Range("I$row")
    $src.Copy()
    $dst.PasteSpecial(-4122)
    $dst.Value = $studentIds[$row]
    $dst.HorizontalAlignment = -4108
    $dst.NumberFormat = "General"
}

$ws.Range("I2:I15").Select()
$excel.ActiveWindow.RangeSelection.Item(1).Activate()
